$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $text) {
    $ws.Range($cell).Value = "'" + $text
    $ws.Range($cell).Style = "Normal"
}

Set-TextCell "D2" '28.476.76'
Set-TextCell "E2" '  +2.38%  '

Set-TextCell "D3" '1.827.23'
Set-TextCell "E3" '  +1.72%  '

Set-TextCell "D4" '1.003'
Set-TextCell "E4" '  +0.20%  '

Set-TextCell "D5" '315.48'
Set-TextCell "E5" '  -0.04%  '

Set-TextCell "D6" '1.002'
Set-TextCell "E6" '  +0.07%  '

Set-TextCell "D7" '0.5021'
Set-TextCell "E7" '  -6.36%  '

Set-TextCell "D8" '0.3911'

Set-TextCell "D9" '0.07706'
Set-TextCell "E9" '  +3.44%  '

Set-TextCell "D10" '41.95'
Set-TextCell "E10" '  +1.12%  '

Set-TextCell "D11" '1.111'
Set-TextCell "E11" '  +2.25%  '

Set-TextCell "D12" '20.99'
Set-TextCell "E12" '  +3.05%  '

Set-TextCell "B13" 'BinanceUSD'
Set-TextCell "C13" 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextCell "D13" '1.003'
Set-TextCell "E13" '  +0.25%  '

Set-TextCell "B14" 'Polkadot'
Set-TextCell "C14" 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextCell "D14" '6.255'
Set-TextCell "E14" '  +0.11%  '

Set-TextCell "D15" '7.551'
Set-TextCell "E15" '  +1.10%  '

Set-TextCell "D16" '1.825.12'
Set-TextCell "E16" '  +1.96%  '

Set-TextCell "D17" '93.40'
Set-TextCell "E17" '  +5.60%  '

Set-TextCell "D18" '0.00001080'
Set-TextCell "E18" '  +1.90%  '

Set-TextCell "D19" '0.06605'
Set-TextCell "E19" '  +1.23%  '

Set-TextCell "D20" '17.69'
Set-TextCell "E20" '  +1.90%  '

Set-TextCell "D21" '1.002'
Set-TextCell "E21" '  +0.03%  '

Set-TextCell "D22" '6.095'
Set-TextCell "E22" '  +2.04%  '

Set-TextCell "D23" '28.503.59'
Set-TextCell "E23" '  +2.35%  '

Set-TextCell "D24" '11.10'
Set-TextCell "E24" '  -0.13%  '

Set-TextCell "D25" '2.259'
Set-TextCell "E25" '  +8.00%  '

Set-TextCell "D26" '157.09'
Set-TextCell "E26" '  +0.30%  '

Set-TextCell "B27" 'EthereumClassic'
Set-TextCell "C27" 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextCell "D27" '20.58'
Set-TextCell "E27" '  +1.42%  '

Set-TextCell "B28" 'WrappedliquidstakedEther2.0'
Set-TextCell "C28" 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextCell "D28" '2.036.08'
Set-TextCell "E28" '  +1.83%  '

Set-TextCell "B29" 'LidoDAOToken'
Set-TextCell "C29" 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextCell "D29" '2.429'
Set-TextCell "E29" '  +4.08%  '

Set-TextCell "D30" '124.70'
Set-TextCell "E30" '  +2.47%  '

Set-TextCell "D31" '1.127'
Set-TextCell "E31" '  +0.80%  '

Set-TextCell "D32" '0.1089'
Set-TextCell "E32" '  -0.37%  '

Set-TextCell "D33" '5.643'
Set-TextCell "E33" '  +2.08%  '

Set-TextCell "D34" '3.662'
Set-TextCell "E34" '  +0.22%  '

Set-TextCell "D35" '0.07077'
Set-TextCell "E35" '  +0.87%  '

Set-TextCell "D36" '0.2219'
Set-TextCell "E36" '  +0.91%  '

Set-TextCell "D37" '9.035'
Set-TextCell "E37" '  +6.99%  '

Set-TextCell "D38" '0.02319'
Set-TextCell "E38" '  +1.90%  '

Set-TextCell "D39" '5.125'
Set-TextCell "E39" '  +0.80%  '

Set-TextCell "D40" '0.6223'
Set-TextCell "E40" '  +1.78%  '

Set-TextCell "D41" '11.19'
Set-TextCell "E41" '  -1.87%  '

Set-TextCell "D42" '1.189'
Set-TextCell "E42" '  +2.35%  '

Set-TextCell "D43" '1.002'

Set-TextCell "D44" '1.398'
Set-TextCell "E44" '  -0.99%  '

Set-TextCell "D45" '13.37'
Set-TextCell "E45" '  +0.78%  '

Set-TextCell "B46" 'PancakeSwap'
Set-TextCell "C46" 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextCell "D46" '3.717'
Set-TextCell "E46" '  +1.06%  '

Set-TextCell "B47" 'Decentraland'
Set-TextCell "C47" 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextCell "D47" '0.5880'
Set-TextCell "E47" '  +2.84%  '

Set-TextCell "D48" '124.26'
Set-TextCell "E48" '  -0.65%  '

Set-TextCell "D49" '1.969'
Set-TextCell "E49" '  +2.98%  '

Set-TextCell "D50" '1.181'
Set-TextCell "E50" '  +0.80%  '

Set-TextCell "D51" '0.06924'
Set-TextCell "E51" '  +1.95%  '

